# Add NLP resources to excel, ivan
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo in existing row 2 "Main Idea" cell (Enbedding -> Embedding) ---
$ws.Range("B2").Value = "1. 这篇文章可以理解Embedding、Attention做了什么事情。`n2."

# --- New row 5: BiDAF code collection (paperswithcode, anchored at #code) ---
$ws.Cells.Item(5, 1).Value = "https://paperswithcode.com/paper/bidirectional-attention-flow-for-machine#code"
$ws.Cells.Item(5, 2).Value = "BiDAF Github代码集合"
$ws.Cells.Item(5, 4).Value = "github"
$ws.Hyperlinks.Add($ws.Range("A5"), "https://paperswithcode.com/paper/bidirectional-attention-flow-for-machine", "code")
$ws.Range("A5:E5").RowHeight = 31.2

# --- New row 6: zhihu attention article ---
$ws.Cells.Item(6, 1).Value = "https://zhuanlan.zhihu.com/p/37601161"
$ws.Cells.Item(6, 2).Value = "深度学习中的注意力模型"
$ws.Cells.Item(6, 3).Value = "张俊林"
$ws.Cells.Item(6, 4).Value = "知乎"
$ws.Cells.Item(6, 5).Value = 201
$ws.Hyperlinks.Add($ws.Range("A6"), "https://zhuanlan.zhihu.com/p/37601161")

# --- New row 7: zhihu word embedding article ---
$ws.Cells.Item(7, 1).Value = "https://zhuanlan.zhihu.com/p/49271699"
$ws.Cells.Item(7, 2).Value = "word embedding技术的发展历史"
$ws.Cells.Item(7, 3).Value = "张俊林"
$ws.Cells.Item(7, 4).Value = "知乎"
$ws.Hyperlinks.Add($ws.Range("A7"), "https://zhuanlan.zhihu.com/p/49271699")

# --- New row 8: Judit Acs masking attention blog post ---
$ws.Cells.Item(8, 1).Value = "http://juditacs.github.io/2018/12/27/masked-attention.html"
$ws.Cells.Item(8, 2).Value = "Masking attention weights in PyTorch"
$ws.Cells.Item(8, 3).Value = "Judit  Acs"
$ws.Cells.Item(8, 4).Value = "github io"
$ws.Cells.Item(8, 5).Value = 2018
$ws.Hyperlinks.Add($ws.Range("A8"), "http://juditacs.github.io/2018/12/27/masked-attention.html")
$ws.Range("A8:E8").RowHeight = 31.2

# --- Page setup (printer defaults that now show up in the saved file) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection cursor left where the author last clicked ---
$ws.Range("C9").Select()
